# to-dolist.xlsx : add weekend notes for 9-10 Sept 2023 (stt 5 & 6 + two
# tip-and-trick entries), per commit "Note 2 ngay cuoi tuan 9-10 Sept 2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 2 : English task now has an "Overview Target" of OK -------------
# --- row 5 : booking task gets its task/overview filled in ----------------
# --- row 6..8 : brand new rows for the two weekend days -------------------
#
# Cells are written in the same order the original author must have typed
# them in (this is what decides the order new entries land in
# xl/sharedStrings.xml), so we reproduce that order here too.

$ws.Range("D5").Value = "Tóm tắt uốn sách lối sống tối giản của Steve Job"
$ws.Range("D6").Value = "Bai hát chari chari lady"
$ws.Range("E5").Value = "Ok"
$ws.Range("D7").Value = "The 5 magical apps that changed my life`nVideo youtube"
$ws.Range("E2").Value = "OK"
$ws.Range("D8").Value = "7 level of reading book.`nHow to remember every thing what you read"

# --- rest of the plain cell values for the two new rows -------------------
$ws.Range("A6").Value = 5
$ws.Range("B6").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B6").Value = 45179
$ws.Range("C6").Value = "English"
$ws.Range("E6").Value = "Ok"

$ws.Range("A7").Value = 6
$ws.Range("B7").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B7").Value = 45179
$ws.Range("C7").Value = "Tip and trick"

$ws.Range("B8").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B8").Value = 45179
$ws.Range("C8").Value = "Tip and trick"

# --- rich-text colour run inside D7 ("Video youtube" in red) --------------
$chars = $ws.Range("D7").Characters(41, 13)
$chars.Font.Color = 255

# --- wrap text on the two long note cells ----------------------------------
$ws.Range("D7").WrapText = $true
$ws.Range("D8").WrapText = $true

# --- row heights for the wrapped rows --------------------------------------
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 30

# --- column widths (new stt/time column layout) -----------------------------
$ws.Columns.Item(1).ColumnWidth = 7.666666666666667
$ws.Columns.Item(2).ColumnWidth = 8.833333333333334
$ws.Columns.Item(3).ColumnWidth = 15.666666666666666
$ws.Columns.Item(4).ColumnWidth = 42.166666666666664
$ws.Columns.Item(5).ColumnWidth = 13.666666666666666

# --- sheet view / selection -------------------------------------------------
[void]($ws.Range("D9").Select())

# --- page setup --------------------------------------------------------------
$ws.PageSetup.Orientation = 1

Write-Output "done"
